$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1421.7142
$ws.Range("I62").Value = 1285.409
$ws.Range("J62").Value = 1921.5
$ws.Range("K62").Value = 1285.409
$ws.Range("L62").Value = 1921.5
$ws.Range("M62").Value = -661.4090000000001
$ws.Range("N62").Value = -3169.5
$ws.Range("H65").Value = 1421.7142
$ws.Range("I65").Value = 1285.409
$ws.Range("J65").Value = 1921.5
$ws.Range("K65").Value = 6427.045
$ws.Range("L65").Value = 9607.5
$ws.Range("M65").Value = -3307.045
$ws.Range("N65").Value = -15847.5
$ws.Range("H98").Value = 62123.332
$ws.Range("I98").Value = 88013
$ws.Range("K98").Value = 88013
$ws.Range("M98").Value = -86515
$ws.Range("H121").Value = 1616.5555
$ws.Range("J121").Value = 1616.5555
$ws.Range("L121").Value = 4849.666499999999
$ws.Range("N121").Value = -8343.666499999999
$ws.Range("H122").Value = 62123.332
$ws.Range("I122").Value = 88013
$ws.Range("K122").Value = 264039
$ws.Range("M122").Value = -261589
$ws.Range("H132").Value = 1833101.5
$ws.Range("I132").Value = 2748006.2
$ws.Range("J132").Value = 3292.077
$ws.Range("K132").Value = 8244018.600000001
$ws.Range("L132").Value = 9876.231
$ws.Range("M132").Value = -8241488.600000001
$ws.Range("N132").Value = -14936.231
$ws.Range("H138").Value = 1952.3906
$ws.Range("I138").Value = 1146.1562
$ws.Range("J138").Value = 2758.625
$ws.Range("K138").Value = 3438.4686
$ws.Range("L138").Value = 8275.875
$ws.Range("M138").Value = 1701.5314
$ws.Range("N138").Value = -18555.875
$ws.Range("H141").Value = 24696274
$ws.Range("I141").Value = 38650624
$ws.Range("J141").Value = 7808.846
$ws.Range("K141").Value = 115951872
$ws.Range("L141").Value = 23426.538
$ws.Range("M141").Value = -115946692
$ws.Range("N141").Value = -33786.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12124.931
$ws.Range("I32").Value = 6543.5225
$ws.Range("J32").Value = 31806.736
$ws.Range("K32").Value = 6543.5225
$ws.Range("L32").Value = 31806.736
$ws.Range("M32").Value = -6256.5225
$ws.Range("N32").Value = -32380.736
$ws.Range("H74").Value = 23811070
$ws.Range("I74").Value = 27779138
$ws.Range("J74").Value = 2664.6667
$ws.Range("K74").Value = 27779138
$ws.Range("L74").Value = 2664.6667
$ws.Range("M74").Value = -27778264
$ws.Range("N74").Value = -4412.6667
$ws.Range("H77").Value = 23811070
$ws.Range("I77").Value = 27779138
$ws.Range("J77").Value = 2664.6667
$ws.Range("K77").Value = 138895690
$ws.Range("L77").Value = 13323.3335
$ws.Range("M77").Value = -138891322
$ws.Range("N77").Value = -22059.3335
$ws.Range("H97").Value = 1092.1072
$ws.Range("I97").Value = 916.0952
$ws.Range("J97").Value = 1620.1428
$ws.Range("K97").Value = 916.0952
$ws.Range("L97").Value = 1620.1428
$ws.Range("M97").Value = -420.0952
$ws.Range("N97").Value = -2612.1428
$ws.Range("H132").Value = 1859.7307
$ws.Range("I132").Value = 1084.6578
$ws.Range("K132").Value = 3253.9734
$ws.Range("M132").Value = -723.9733999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2275144.2
$ws.Range("I105").Value = 3789655.2
$ws.Range("J105").Value = 3377.75
$ws.Range("K105").Value = 3789655.2
$ws.Range("L105").Value = 3377.75
$ws.Range("M105").Value = -3787908.2
$ws.Range("N105").Value = -6871.75
$ws.Range("H134").Value = 8335253
$ws.Range("I134").Value = 10418216
$ws.Range("J134").Value = 3400
$ws.Range("K134").Value = 31254648
$ws.Range("L134").Value = 10200
$ws.Range("M134").Value = -31252113
$ws.Range("N134").Value = -15270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1758333.4
$ws.Range("J4").Value = 1827272.8
$ws.Range("L4").Value = 1827272.8
$ws.Range("N4").Value = -1827496.8
$ws.Range("H107").Value = 1103.375
$ws.Range("I107").Value = 583.1539
$ws.Range("J107").Value = 1718.1818
$ws.Range("K107").Value = 583.1539
$ws.Range("L107").Value = 1718.1818
$ws.Range("M107").Value = 1336.8461
$ws.Range("N107").Value = -5558.1818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 174.1
$ws.Range("I4").Value = 174.1
$ws.Range("K4").Value = 522.3
$ws.Range("M4").Value = -410.3
$ws.Range("H5").Value = 472.77142
$ws.Range("I5").Value = 298.69565
$ws.Range("J5").Value = 806.4167
$ws.Range("K5").Value = 896.08695
$ws.Range("L5").Value = 2419.2501
$ws.Range("M5").Value = -784.08695
$ws.Range("N5").Value = -2643.2501
$ws.Range("H57").Value = 5541.7896
$ws.Range("I57").Value = 1002.5
$ws.Range("J57").Value = 6075.8237
$ws.Range("K57").Value = 3007.5
$ws.Range("L57").Value = 18227.4711
$ws.Range("M57").Value = -2448.5
$ws.Range("N57").Value = -19345.4711
$ws.Range("H68").Value = 5597.421
$ws.Range("I68").Value = 8516.833000000001
$ws.Range("J68").Value = 592.7143
$ws.Range("K68").Value = 25550.499
$ws.Range("L68").Value = 1778.1429
$ws.Range("M68").Value = -24739.499
$ws.Range("N68").Value = -3400.1429
$ws.Range("H71").Value = 5597.421
$ws.Range("I71").Value = 8516.833000000001
$ws.Range("J71").Value = 592.7143
$ws.Range("K71").Value = 76651.497
$ws.Range("L71").Value = 5334.428699999999
$ws.Range("M71").Value = -72595.497
$ws.Range("N71").Value = -13446.4287
$ws.Range("H113").Value = 23809956
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 35714660
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 107143980
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -107148320
$ws.Range("H131").Value = 13945195
$ws.Range("I131").Value = 50000184
$ws.Range("J131").Value = 77890.46000000001
$ws.Range("K131").Value = 150000552
$ws.Range("L131").Value = 233671.38
$ws.Range("M131").Value = -149995512
$ws.Range("N131").Value = -243751.38
$ws.Range("H135").Value = 472.77142
$ws.Range("I135").Value = 298.69565
$ws.Range("J135").Value = 806.4167
$ws.Range("K135").Value = 2688.26085
$ws.Range("L135").Value = 7257.7503
$ws.Range("M135").Value = -153.2608500000001
$ws.Range("N135").Value = -12327.7503
$ws.Range("H139").Value = 2178.5476
$ws.Range("I139").Value = 1374.2142
$ws.Range("J139").Value = 3787.2144
$ws.Range("K139").Value = 4122.642599999999
$ws.Range("L139").Value = 11361.6432
$ws.Range("M139").Value = 1017.357400000001
$ws.Range("N139").Value = -21641.6432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8080
$ws.Range("J5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("N5").Value = -10224
$ws.Range("H70").Value = 4119.3
$ws.Range("I70").Value = 4199.7
$ws.Range("J70").Value = 4038.9
$ws.Range("K70").Value = 4199.7
$ws.Range("L70").Value = 4038.9
$ws.Range("M70").Value = -3929.7
$ws.Range("N70").Value = -4578.9
$ws.Range("H73").Value = 4119.3
$ws.Range("I73").Value = 4199.7
$ws.Range("J73").Value = 4038.9
$ws.Range("K73").Value = 4199.7
$ws.Range("L73").Value = 4038.9
$ws.Range("M73").Value = -3263.7
$ws.Range("N73").Value = -5910.9
$ws.Range("H113").Value = 1404.9697
$ws.Range("I113").Value = 1621.8125
$ws.Range("J113").Value = 1200.8823
$ws.Range("K113").Value = 1621.8125
$ws.Range("L113").Value = 1200.8823
$ws.Range("M113").Value = 548.1875
$ws.Range("N113").Value = -5540.8823
$ws.Range("H126").Value = 2121.65
$ws.Range("I126").Value = 1802
$ws.Range("J126").Value = 2441.3
$ws.Range("K126").Value = 5406
$ws.Range("L126").Value = 7323.900000000001
$ws.Range("M126").Value = -2936
$ws.Range("N126").Value = -12263.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10916.667
$ws.Range("J2").Value = 10916.667
$ws.Range("L2").Value = 10916.667
$ws.Range("N2").Value = -11140.667
$ws.Range("H55").Value = 420.27274
$ws.Range("I55").Value = 355.125
$ws.Range("K55").Value = 355.125
$ws.Range("M55").Value = -182.125
$ws.Range("H132").Value = 2751.6
$ws.Range("I132").Value = 1987.1765
$ws.Range("K132").Value = 5961.529500000001
$ws.Range("M132").Value = -3431.529500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5004995
$ws.Range("J2").Value = 10000000
$ws.Range("L2").Value = 10000000
$ws.Range("N2").Value = -10000224
$ws.Range("H122").Value = 1463.909
$ws.Range("I122").Value = 1137.875
$ws.Range("K122").Value = 3413.625
$ws.Range("M122").Value = -963.625
$ws.Range("H132").Value = 1422.475
$ws.Range("I132").Value = 883.6129
$ws.Range("J132").Value = 3278.5557
$ws.Range("K132").Value = 2650.8387
$ws.Range("L132").Value = 9835.667099999999
$ws.Range("M132").Value = -120.8386999999998
$ws.Range("N132").Value = -14895.6671
